$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 160
$ws.Range("I19").Value = 80
$ws.Range("K19").Value = 80
$ws.Range("M19").Value = 95

$ws.Range("H38").Value = 250.9
$ws.Range("I38").Value = 178.77777
$ws.Range("J38").Value = 900
$ws.Range("K38").Value = 536.33331
$ws.Range("L38").Value = 2700
$ws.Range("M38").Value = -164.33331
$ws.Range("N38").Value = -3444

$ws.Range("H46").Value = 1150
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1150
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3450
$ws.Range("N46").Value = -3688
$ws.Range("M46").ClearContents()

$ws.Range("H60").Value = 1150
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 1150
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 3450
$ws.Range("N60").Value = -4418
$ws.Range("M60").ClearContents()

$ws.Range("H100").Value = 2552.9412
$ws.Range("I100").Value = 1860
$ws.Range("K100").Value = 1860
$ws.Range("M100").Value = -1319

$ws.Range("H106").Value = 1543.7693
$ws.Range("I106").Value = 1156.4
$ws.Range("K106").Value = 1156.4
$ws.Range("M106").Value = -525.4000000000001

$ws.Range("H113").Value = 90913190
$ws.Range("I113").Value = 142859540
$ws.Range("J113").Value = 7097
$ws.Range("K113").Value = 142859540
$ws.Range("L113").Value = 7097
$ws.Range("M113").Value = -142856286
$ws.Range("N113").Value = -13605

$ws.Range("H116").Value = 3447.111
$ws.Range("I116").Value = 2004
$ws.Range("K116").Value = 2004
$ws.Range("M116").Value = 1438

$ws.Range("H129").Value = 701.375
$ws.Range("J129").Value = 875
$ws.Range("L129").Value = 2625
$ws.Range("N129").Value = -12625

$ws.Range("H132").Value = 59604.61
$ws.Range("I132").Value = 59604.61
$ws.Range("K132").Value = 178813.83
$ws.Range("M132").Value = -176283.83

$ws.Range("H138").Value = 2590.0256
$ws.Range("J138").Value = 3007.7188
$ws.Range("L138").Value = 9023.1564
$ws.Range("N138").Value = -19303.1564

$ws.Range("H141").Value = 1700.7894
$ws.Range("I141").Value = 1077.3529
$ws.Range("K141").Value = 3232.0587
$ws.Range("M141").Value = 1947.9413

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H74").Value = 625.27026
$ws.Range("I74").Value = 342.32144
$ws.Range("J74").Value = 1505.5555
$ws.Range("K74").Value = 342.32144
$ws.Range("L74").Value = 1505.5555
$ws.Range("M74").Value = 531.6785600000001
$ws.Range("N74").Value = -3253.5555

$ws.Range("H77").Value = 625.27026
$ws.Range("I77").Value = 342.32144
$ws.Range("J77").Value = 1505.5555
$ws.Range("K77").Value = 1711.6072
$ws.Range("L77").Value = 7527.7775
$ws.Range("M77").Value = 2656.3928
$ws.Range("N77").Value = -16263.7775

$ws.Range("H132").Value = 27168.55
$ws.Range("I132").Value = 1699.5
$ws.Range("J132").Value = 65372.125
$ws.Range("K132").Value = 5098.5
$ws.Range("L132").Value = 196116.375
$ws.Range("M132").Value = -2568.5
$ws.Range("N132").Value = -201176.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 50000
$ws.Range("J27").Value = 50000
$ws.Range("L27").Value = 50000
$ws.Range("N27").Value = -50384

$ws.Range("H94").Value = 4820.643
$ws.Range("I94").Value = 1999.75
$ws.Range("J94").Value = 5949
$ws.Range("K94").Value = 1999.75
$ws.Range("L94").Value = 5949
$ws.Range("M94").Value = -1548.75
$ws.Range("N94").Value = -6851

$ws.Range("H134").Value = 19387.492
$ws.Range("I134").Value = 21118.21
$ws.Range("K134").Value = 63354.63
$ws.Range("M134").Value = -60819.63

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10002.214
$ws.Range("I31").Value = 12890.741
$ws.Range("K31").Value = 12890.741
$ws.Range("M31").Value = -12595.741

$ws.Range("H34").Value = 10002.214
$ws.Range("I34").Value = 12890.741
$ws.Range("K34").Value = 12890.741
$ws.Range("M34").Value = -12688.741

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 4760.7085
$ws.Range("J107").Value = 724.7273
$ws.Range("L107").Value = 2174.1819
$ws.Range("N107").Value = -6014.1819

$ws.Range("H109").Value = 6077
$ws.Range("J109").Value = 6499.125
$ws.Range("L109").Value = 19497.375
$ws.Range("N109").Value = -21577.375

$ws.Range("H131").Value = 157096.12
$ws.Range("J131").Value = 170326.42
$ws.Range("L131").Value = 510979.26
$ws.Range("N131").Value = -521059.26

$ws.Range("H138").Value = 1559
$ws.Range("I138").Value = 1559
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 4677
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 463
$ws.Range("N138").ClearContents()

$ws.Range("H140").Value = 5385.5
$ws.Range("I140").Value = 6889.3125
$ws.Range("K140").Value = 20667.9375
$ws.Range("M140").Value = -15487.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4404.4287
$ws.Range("I70").Value = 3961
$ws.Range("J70").Value = 4581.8
$ws.Range("K70").Value = 3961
$ws.Range("L70").Value = 4581.8
$ws.Range("M70").Value = -3691
$ws.Range("N70").Value = -5121.8

$ws.Range("H73").Value = 4404.4287
$ws.Range("I73").Value = 3961
$ws.Range("J73").Value = 4581.8
$ws.Range("K73").Value = 3961
$ws.Range("L73").Value = 4581.8
$ws.Range("M73").Value = -3025
$ws.Range("N73").Value = -6453.8

$ws.Range("H102").Value = 2606.1333
$ws.Range("I102").Value = 2915.9048
$ws.Range("J102").Value = 1883.3334
$ws.Range("K102").Value = 2915.9048
$ws.Range("L102").Value = 1883.3334
$ws.Range("M102").Value = -1293.9048
$ws.Range("N102").Value = -5127.3334

$ws.Range("H126").Value = 4754.1665
$ws.Range("I126").Value = 3640
$ws.Range("J126").Value = 6611.1113
$ws.Range("K126").Value = 10920
$ws.Range("L126").Value = 19833.3339
$ws.Range("M126").Value = -8450
$ws.Range("N126").Value = -24773.3339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 272.55554
$ws.Range("I55").Value = 201.22223
$ws.Range("J55").Value = 343.8889
$ws.Range("K55").Value = 201.22223
$ws.Range("L55").Value = 343.8889
$ws.Range("M55").Value = -28.22223
$ws.Range("N55").Value = -689.8888999999999

$ws.Range("H93").Value = 3450
$ws.Range("I93").Value = 3570.8572
$ws.Range("J93").Value = 3280.8
$ws.Range("K93").Value = 3570.8572
$ws.Range("L93").Value = 3280.8
$ws.Range("M93").Value = -2322.8572
$ws.Range("N93").Value = -5776.8

$ws.Range("H122").Value = 2634.9565
$ws.Range("I122").Value = 2280.4
$ws.Range("J122").Value = 2907.6924
$ws.Range("K122").Value = 6841.200000000001
$ws.Range("L122").Value = 8723.0772
$ws.Range("M122").Value = -4391.200000000001
$ws.Range("N122").Value = -13623.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5169.75
$ws.Range("J62").Value = 4899.5
$ws.Range("L62").Value = 4899.5
$ws.Range("N62").Value = -6147.5

$ws.Range("H65").Value = 5169.75
$ws.Range("J65").Value = 4899.5
$ws.Range("L65").Value = 24497.5
$ws.Range("N65").Value = -30737.5

$ws.Range("H96").Value = 4700
$ws.Range("I96").Value = 2250
$ws.Range("J96").Value = 5312.5
$ws.Range("K96").Value = 2250
$ws.Range("L96").Value = 5312.5
$ws.Range("M96").Value = -877
$ws.Range("N96").Value = -8058.5
